$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")
$ws.Activate()

# Update data row (A2:C2) first, then header (D1), then new data cell (D2),
# matching the order the shared strings were authored in the source workbook.
$ws.Range("A2").Value = "Houlihan Capital Holdings, Inc"
$ws.Range("B2").Value = "Houlihan"
$ws.Range("C2").Value = "Employee"
$ws.Range("D1").Value = "ContactName"
$ws.Range("D2").Value = "Houlihan Employee"

# Make the new D1 header bold to match the rest of the header row
$ws.Range("D1").Font.Bold = $true

# Adjust column widths to fit the new, wider content (bestFit/AutoFit-style resize)
$ws.Columns.Item(1).ColumnWidth = 25.44140625
$ws.Columns.Item(4).ColumnWidth = 16.77734375

# Update the active selection shown when the sheet is viewed
$ws.Range("D7").Select()

$wb.Save()
